$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26 (shifts rows 26.. down by one),
# copying formatting from the row above (row 25), like Excel's default behavior.
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new component data.
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = "100kΩ±5%   "
$ws.Cells.Item(26, 6).Value = "YAGEO"
$ws.Cells.Item(26, 7).Value = "0805"
$ws.Cells.Item(26, 8).Value = "RC0805JR-07100KL"

# Update the selected cell, matching the author's final selection.
$ws.Range("I35").Select()
